# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.613.16'
$ws.Range("E2").Value = '  +3.84%  '
$ws.Range("D3").Value = '3.791.23'
$ws.Range("E3").Value = '  +7.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '420.07'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.50'
$ws.Range("E6").Value = '  +4.23%  '
$ws.Range("D7").Value = '3.775.42'
$ws.Range("E7").Value = '  +7.27%  '
$ws.Range("E8").Value = '  -0.97%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  +13.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000412'
$ws.Range("E12").Value = '  +51.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.24'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.51'
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("D15").Value = '4.384.54'
$ws.Range("E15").Value = '  +7.27%  '
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").Value = '3.756.62'
$ws.Range("E17").Value = '  +6.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.54'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.12'
$ws.Range("E19").Value = '  +2.75%  '
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").Value = '68.447.28'
$ws.Range("E21").Value = '  +3.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '446.23'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.51'
$ws.Range("E23").Value = '  +17.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '90.64'
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.10'
$ws.Range("E25").Value = '  -4.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.45'
$ws.Range("E26").Value = '  +12.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.32'
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.10'
$ws.Range("E29").Value = '  +5.46%  '
$ws.Range("E30").Value = '  +5.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.64'
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.15'
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.164'
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.39'
$ws.Range("E35").Value = '  +5.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.07'
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0490'
$ws.Range("E38").Value = '  -2.82%  '
$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("E39").Value = '  +30.41%  '
$ws.Range("D40").Value = '0.0₃0715'
$ws.Range("E40").Value = '  -3.66%  '
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.995'
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("B43").Value = 'LidoDAOToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.38'
$ws.Range("E43").Value = '  +3.97%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.39'
$ws.Range("E44").Value = '  +27.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.22'
$ws.Range("E45").Value = '  +1.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.19'
$ws.Range("E46").Value = '  +24.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.90'
$ws.Range("E47").Value = '  -4.44%  '
$ws.Range("E48").Value = '  +4.75%  '
$ws.Range("E49").Value = '  -3.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.62'
$ws.Range("E50").Value = '  -6.44%  '
$ws.Range("E51").Value = '  -2.23%  '
